$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.907.58'

$ws.Range("D3").Value = '2.233.51'
$ws.Range("E3").Value = '  -0.55%  '

$ws.Range("E4").Value = '  -0.24%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '232.50'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.47%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.625'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -1.77%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '60.62'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -6.78%  '

$ws.Range("E8").Value = '  -0.07%  '

$ws.Range("E9").Value = '  -0.30%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '58.17'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -4.36%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0903'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -0.31%  '

$ws.Range("E12").Value = '  -0.61%  '

$ws.Range("D13").Value = '2.565.53'
$ws.Range("E13").Value = '  -0.38%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '15.55'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -3.70%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '22.74'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +1.85%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.65'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +0.44%  '

$ws.Range("D18").Value = '2.242.33'
$ws.Range("E18").Value = '  -0.09%  '

$ws.Range("D19").Value = '41.891.50'
$ws.Range("E19").Value = '  +1.33%  '

$ws.Range("D20").Value = '0.0₃0914'
$ws.Range("E20").Value = '  -0.27%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.70'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -1.79%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.18'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.59%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '248.74'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -1.80%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.999'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -0.21%  '

$ws.Range("B25").Value = 'Toncoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.40'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +2.91%  '

$ws.Range("B26").Value = 'PancakeSwap'
$ws.Range("C26").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.39'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -0.06%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.68'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -0.80%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '169.33'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -1.97%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.141'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -2.56%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '19.96'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -2.46%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.41'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -2.18%  '

$ws.Range("E32").Value = '  -10.90%  '

$ws.Range("E33").Value = '  -1.58%  '

$ws.Range("E34").Value = '  +3.57%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.72'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +0.32%  '

$ws.Range("E36").Value = '  +4.39%  '

$ws.Range("E37").Value = '  -8.60%  '

$ws.Range("E38").Value = '  -1.69%  '

$ws.Range("E39").Value = '  -5.22%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.000247'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +5.56%  '

$ws.Range("E41").Value = '  -0.20%  '

$ws.Range("E42").Value = '  +2.05%  '

$ws.Range("E43").Value = '  -0.49%  '

$ws.Range("B44").Value = 'TrustWalletToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.23'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -0.20%  '

$ws.Range("B45").Value = 'FTXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.52'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -6.68%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '99.13'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -2.68%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0969'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +2.84%  '

$ws.Range("D48").Value = '1.474.61'
$ws.Range("E48").Value = '  -2.43%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '16.71'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -6.15%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.29'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +8.68%  '

$ws.Range("E51").Value = '  -2.60%  '
